$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Modify existing rows (DUPE flags / corrections) ---
$ws.Cells.Item(1090, 8).Value = "Yes"   # H1090: No -> Yes

$ws.Cells.Item(1149, 4).Value = "Pontiac"  # D1149: Fake Generic Brand -> Pontiac
$ws.Cells.Item(1149, 6).Value = 1991       # F1149: 2015 -> 1991
$ws.Cells.Item(1149, 8).Value = "Yes"      # H1149: No -> Yes

# --- Append new rows 1300-1341 ---
# Row 1300
$ws.Cells.Item(1300, 1).Value = "HotWheels"
$ws.Cells.Item(1300, 2).Value = 2021
$ws.Cells.Item(1300, 3).Value = "Mainline (M Case)"
$ws.Cells.Item(1300, 4).Value = "Ford"
$ws.Cells.Item(1300, 5).Value = "Bronco"
$ws.Cells.Item(1300, 6).Value = 2021
$ws.Cells.Item(1300, 7).Value = "Orange"
$ws.Cells.Item(1300, 8).Value = "No"

# Row 1301
$ws.Cells.Item(1301, 1).Value = "HotWheels"
$ws.Cells.Item(1301, 2).Value = 2022
$ws.Cells.Item(1301, 3).Value = "Mainline (N Case)"
$ws.Cells.Item(1301, 4).Value = "Chevy"
$ws.Cells.Item(1301, 5).Value = "Camaro Convertible"
$ws.Cells.Item(1301, 6).Value = 1969
$ws.Cells.Item(1301, 7).Value = "Black"
$ws.Cells.Item(1301, 8).Value = "No"

# Row 1302
$ws.Cells.Item(1302, 1).Value = "HotWheels"
$ws.Cells.Item(1302, 2).Value = 2023
$ws.Cells.Item(1302, 3).Value = "Mainline (Q Case)"
$ws.Cells.Item(1302, 4).Value = "Honda"
$ws.Cells.Item(1302, 5).Value = "Super Cub"
$ws.Cells.Item(1302, 6).Value = 2016
$ws.Cells.Item(1302, 7).Value = "Yellow"
$ws.Cells.Item(1302, 8).Value = "No"

# Row 1303
$ws.Cells.Item(1303, 1).Value = "HotWheels"
$ws.Cells.Item(1303, 2).Value = 2022
$ws.Cells.Item(1303, 3).Value = "Color Shifters"
$ws.Cells.Item(1303, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1303, 5).Value = "Scorpedo"
$ws.Cells.Item(1303, 6).Value = 2022
$ws.Cells.Item(1303, 7).Value = "Orange"
$ws.Cells.Item(1303, 8).Value = "No"

# Row 1304
$ws.Cells.Item(1304, 1).Value = "Johnny Lightning"
$ws.Cells.Item(1304, 2).Value = 2021
$ws.Cells.Item(1304, 3).Value = "OK Used Cars (Rel 4)"
$ws.Cells.Item(1304, 4).Value = "Pontiac"
$ws.Cells.Item(1304, 5).Value = "Firebird T/A WS6"
$ws.Cells.Item(1304, 6).Value = 1997
$ws.Cells.Item(1304, 7).Value = "Black"
$ws.Cells.Item(1304, 8).Value = "No"

# Row 1305
$ws.Cells.Item(1305, 1).Value = "M2"
$ws.Cells.Item(1305, 2).Value = 2022
$ws.Cells.Item(1305, 3).Value = "Machines"
$ws.Cells.Item(1305, 4).Value = "Chevy"
$ws.Cells.Item(1305, 5).Value = "Camaro IROC-Z"
$ws.Cells.Item(1305, 6).Value = 1985
$ws.Cells.Item(1305, 7).Value = "Dark Green Millitary"
$ws.Cells.Item(1305, 8).Value = "No"

# Row 1306
$ws.Cells.Item(1306, 1).Value = "HotWheels"
$ws.Cells.Item(1306, 2).Value = 1997
$ws.Cells.Item(1306, 3).Value = "HotWheels 30 Years"
$ws.Cells.Item(1306, 4).Value = "Ford"
$ws.Cells.Item(1306, 5).Value = "Vicky"
$ws.Cells.Item(1306, 6).Value = 1932
$ws.Cells.Item(1306, 7).Value = "Red"
$ws.Cells.Item(1306, 8).Value = "No"

# Row 1307
$ws.Cells.Item(1307, 1).Value = "HotWheels"
$ws.Cells.Item(1307, 2).Value = 2013
$ws.Cells.Item(1307, 3).Value = "The Jetsons"
$ws.Cells.Item(1307, 4).Value = "Fake Generic Brand"
$ws.Cells.Item(1307, 5).Value = "Capsule Car"
$ws.Cells.Item(1307, 6).Value = 2013
$ws.Cells.Item(1307, 7).Value = "Lime Green"
$ws.Cells.Item(1307, 8).Value = "No"

# Row 1308
$ws.Cells.Item(1308, 1).Value = "HotWheels"
$ws.Cells.Item(1308, 2).Value = 2017
$ws.Cells.Item(1308, 3).Value = "Mainline (B Case)"
$ws.Cells.Item(1308, 4).Value = "Pontiac"
$ws.Cells.Item(1308, 5).Value = "K.I.T.T."
$ws.Cells.Item(1308, 6).Value = 1991
$ws.Cells.Item(1308, 7).Value = "Black"
$ws.Cells.Item(1308, 8).Value = "No"

# Row 1309
$ws.Cells.Item(1309, 1).Value = "HotWheels"
$ws.Cells.Item(1309, 2).Value = 2001
$ws.Cells.Item(1309, 3).Value = "Mainline"
$ws.Cells.Item(1309, 4).Value = "Chevy"
$ws.Cells.Item(1309, 5).Value = "Custom Pickup"
$ws.Cells.Item(1309, 6).Value = 1969
$ws.Cells.Item(1309, 7).Value = "Red"
$ws.Cells.Item(1309, 8).Value = "No"

# Row 1310
$ws.Cells.Item(1310, 1).Value = "HotWheels"
$ws.Cells.Item(1310, 2).Value = 2021
$ws.Cells.Item(1310, 3).Value = "Mainline (M Case)"
$ws.Cells.Item(1310, 4).Value = "Chevy"
$ws.Cells.Item(1310, 5).Value = "C10"
$ws.Cells.Item(1310, 6).Value = 1967
$ws.Cells.Item(1310, 7).Value = "Light Blue"
$ws.Cells.Item(1310, 8).Value = "No"

# Row 1311
$ws.Cells.Item(1311, 1).Value = "HotWheels"
$ws.Cells.Item(1311, 2).Value = 2015
$ws.Cells.Item(1311, 3).Value = "Mainline (M Case)"
$ws.Cells.Item(1311, 4).Value = "Chevy"
$ws.Cells.Item(1311, 5).Value = "Camaro"
$ws.Cells.Item(1311, 6).Value = 1981
$ws.Cells.Item(1311, 7).Value = "Black"
$ws.Cells.Item(1311, 8).Value = "No"

# Row 1312
$ws.Cells.Item(1312, 1).Value = "Takara Tomy"
$ws.Cells.Item(1312, 2).Value = 2016
$ws.Cells.Item(1312, 3).Value = "Main"
$ws.Cells.Item(1312, 4).Value = "Toyota"
$ws.Cells.Item(1312, 5).Value = "L&F Geneo"
$ws.Cells.Item(1312, 6).Value = 2015
$ws.Cells.Item(1312, 7).Value = "Orange"
$ws.Cells.Item(1312, 8).Value = "No"

# Row 1313
$ws.Cells.Item(1313, 1).Value = "HotWheels"
$ws.Cells.Item(1313, 2).Value = 1997
$ws.Cells.Item(1313, 3).Value = "Mainline First Editions"
$ws.Cells.Item(1313, 4).Value = "Jaguar"
$ws.Cells.Item(1313, 5).Value = "D-Type"
$ws.Cells.Item(1313, 6).Value = 1957
$ws.Cells.Item(1313, 7).Value = "Black"
$ws.Cells.Item(1313, 8).Value = "No"
$ws.Cells.Item(1313, 9).Value = "Damaged"

# Row 1314
$ws.Cells.Item(1314, 1).Value = "HotWheels"
$ws.Cells.Item(1314, 2).Value = 2022
$ws.Cells.Item(1314, 3).Value = "Mainline (N Case)"
$ws.Cells.Item(1314, 4).Value = "Chevy"
$ws.Cells.Item(1314, 5).Value = "Camaro Convertible"
$ws.Cells.Item(1314, 6).Value = 1969
$ws.Cells.Item(1314, 7).Value = "Black"
$ws.Cells.Item(1314, 8).Value = "Yes"

# Row 1315
$ws.Cells.Item(1315, 1).Value = "HotWheels"
$ws.Cells.Item(1315, 2).Value = 1999
$ws.Cells.Item(1315, 3).Value = "Mainline First Editions"
$ws.Cells.Item(1315, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1315, 5).Value = "Phaeton"
$ws.Cells.Item(1315, 6).Value = 1999
$ws.Cells.Item(1315, 7).Value = "Greenish Blue"
$ws.Cells.Item(1315, 8).Value = "No"

# Row 1316
$ws.Cells.Item(1316, 1).Value = "HotWheels"
$ws.Cells.Item(1316, 2).Value = 1995
$ws.Cells.Item(1316, 3).Value = "Mainline"
$ws.Cells.Item(1316, 4).Value = "Ford"
$ws.Cells.Item(1316, 5).Value = "LTL"
$ws.Cells.Item(1316, 6).Value = 1995
$ws.Cells.Item(1316, 7).Value = "Silver #11"
$ws.Cells.Item(1316, 8).Value = "No"

# Row 1317
$ws.Cells.Item(1317, 1).Value = "HotWheels"
$ws.Cells.Item(1317, 2).Value = 2018
$ws.Cells.Item(1317, 3).Value = "Mainline (H Case)"
$ws.Cells.Item(1317, 4).Value = "Dodge"
$ws.Cells.Item(1317, 5).Value = "Charger Drift"
$ws.Cells.Item(1317, 6).Value = 2015
$ws.Cells.Item(1317, 7).Value = "White"
$ws.Cells.Item(1317, 8).Value = "No"

# Row 1318
$ws.Cells.Item(1318, 1).Value = "HotWheels"
$ws.Cells.Item(1318, 2).Value = 2018
$ws.Cells.Item(1318, 3).Value = "Mainline (Q Case)"
$ws.Cells.Item(1318, 4).Value = "Chevy"
$ws.Cells.Item(1318, 5).Value = "Classic Nomad"
$ws.Cells.Item(1318, 6).Value = 1955
$ws.Cells.Item(1318, 7).Value = "Orange"
$ws.Cells.Item(1318, 8).Value = "Yes"

# Row 1319
$ws.Cells.Item(1319, 1).Value = "HotWheels"
$ws.Cells.Item(1319, 2).Value = 2019
$ws.Cells.Item(1319, 3).Value = "Mainline (B Case)"
$ws.Cells.Item(1319, 4).Value = "Shelby"
$ws.Cells.Item(1319, 5).Value = "GT-500"
$ws.Cells.Item(1319, 6).Value = 1967
$ws.Cells.Item(1319, 7).Value = "Blue"
$ws.Cells.Item(1319, 8).Value = "No"

# Row 1320
$ws.Cells.Item(1320, 1).Value = "HotWheels"
$ws.Cells.Item(1320, 2).Value = 2018
$ws.Cells.Item(1320, 3).Value = "Mainline (G Case)"
$ws.Cells.Item(1320, 4).Value = "Ford"
$ws.Cells.Item(1320, 5).Value = "Coupe"
$ws.Cells.Item(1320, 6).Value = 1932
$ws.Cells.Item(1320, 7).Value = "Orange"
$ws.Cells.Item(1320, 8).Value = "No"

# Row 1321
$ws.Cells.Item(1321, 1).Value = "HotWheels"
$ws.Cells.Item(1321, 2).Value = 2018
$ws.Cells.Item(1321, 3).Value = "Mainline (Q Case)"
$ws.Cells.Item(1321, 4).Value = "Chevy"
$ws.Cells.Item(1321, 5).Value = "Corvette C7 Z06"
$ws.Cells.Item(1321, 6).Value = 2015
$ws.Cells.Item(1321, 7).Value = "Blue"
$ws.Cells.Item(1321, 8).Value = "No"

# Row 1322
$ws.Cells.Item(1322, 1).Value = "HotWheels"
$ws.Cells.Item(1322, 2).Value = 2018
$ws.Cells.Item(1322, 3).Value = "Mainline (B Case)"
$ws.Cells.Item(1322, 4).Value = "Ford"
$ws.Cells.Item(1322, 5).Value = "Mustang Mach1"
$ws.Cells.Item(1322, 6).Value = 1971
$ws.Cells.Item(1322, 7).Value = "Red"
$ws.Cells.Item(1322, 8).Value = "No"

# Row 1323
$ws.Cells.Item(1323, 1).Value = "HotWheels"
$ws.Cells.Item(1323, 2).Value = 2017
$ws.Cells.Item(1323, 3).Value = "Character Cars"
$ws.Cells.Item(1323, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1323, 5).Value = "Mosasaurus"
$ws.Cells.Item(1323, 6).Value = 2015
$ws.Cells.Item(1323, 7).Value = "Dark Blue"
$ws.Cells.Item(1323, 8).Value = "No"

# Row 1324
$ws.Cells.Item(1324, 1).Value = "HotWheels"
$ws.Cells.Item(1324, 2).Value = 2017
$ws.Cells.Item(1324, 3).Value = "Character Cars"
$ws.Cells.Item(1324, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1324, 5).Value = "Velociraptor Blue"
$ws.Cells.Item(1324, 6).Value = 2015
$ws.Cells.Item(1324, 7).Value = "Gray"
$ws.Cells.Item(1324, 8).Value = "No"

# Row 1325
$ws.Cells.Item(1325, 1).Value = "HotWheels"
$ws.Cells.Item(1325, 2).Value = 2017
$ws.Cells.Item(1325, 3).Value = "Character Cars"
$ws.Cells.Item(1325, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1325, 5).Value = "Stegosaurus"
$ws.Cells.Item(1325, 6).Value = 2015
$ws.Cells.Item(1325, 7).Value = "Dark Green"
$ws.Cells.Item(1325, 8).Value = "No"

# Row 1326
$ws.Cells.Item(1326, 1).Value = "HotWheels"
$ws.Cells.Item(1326, 2).Value = 2017
$ws.Cells.Item(1326, 3).Value = "Character Cars"
$ws.Cells.Item(1326, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1326, 5).Value = "T-Rex"
$ws.Cells.Item(1326, 6).Value = 2015
$ws.Cells.Item(1326, 7).Value = "Brown"
$ws.Cells.Item(1326, 8).Value = "No"

# Row 1327
$ws.Cells.Item(1327, 1).Value = "HotWheels"
$ws.Cells.Item(1327, 2).Value = 2017
$ws.Cells.Item(1327, 3).Value = "Character Cars"
$ws.Cells.Item(1327, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1327, 5).Value = "Triceratops"
$ws.Cells.Item(1327, 6).Value = 2015
$ws.Cells.Item(1327, 7).Value = "Dark Green"
$ws.Cells.Item(1327, 8).Value = "No"

# Row 1328
$ws.Cells.Item(1328, 1).Value = "HotWheels"
$ws.Cells.Item(1328, 2).Value = 2018
$ws.Cells.Item(1328, 3).Value = "Chevy Trucks 100th"
$ws.Cells.Item(1328, 4).Value = "Chevy"
$ws.Cells.Item(1328, 5).Value = "Silverado"
$ws.Cells.Item(1328, 6).Value = 1983
$ws.Cells.Item(1328, 7).Value = "Black"
$ws.Cells.Item(1328, 8).Value = "No"

# Row 1329
$ws.Cells.Item(1329, 1).Value = "HotWheels"
$ws.Cells.Item(1329, 2).Value = 2018
$ws.Cells.Item(1329, 3).Value = "Chevy Trucks 100th"
$ws.Cells.Item(1329, 4).Value = "Chevy"
$ws.Cells.Item(1329, 5).Value = "Custom Pickup"
$ws.Cells.Item(1329, 6).Value = 1969
$ws.Cells.Item(1329, 7).Value = "Orange"
$ws.Cells.Item(1329, 8).Value = "No"

# Row 1330
$ws.Cells.Item(1330, 1).Value = "HotWheels"
$ws.Cells.Item(1330, 2).Value = 2015
$ws.Cells.Item(1330, 3).Value = "Mainline (A Case)"
$ws.Cells.Item(1330, 4).Value = "Chevy"
$ws.Cells.Item(1330, 5).Value = "Corvette Racer"
$ws.Cells.Item(1330, 6).Value = 1969
$ws.Cells.Item(1330, 7).Value = "Blue #69"
$ws.Cells.Item(1330, 8).Value = "No"

# Row 1331
$ws.Cells.Item(1331, 1).Value = "HotWheels"
$ws.Cells.Item(1331, 2).Value = 2018
$ws.Cells.Item(1331, 3).Value = "Mainline (B Case)"
$ws.Cells.Item(1331, 4).Value = "Fake Generic Brand"
$ws.Cells.Item(1331, 5).Value = "Milano"
$ws.Cells.Item(1331, 6).Value = 2018
$ws.Cells.Item(1331, 7).Value = "Gray"
$ws.Cells.Item(1331, 8).Value = "No"

# Row 1332
$ws.Cells.Item(1332, 1).Value = "HotWheels"
$ws.Cells.Item(1332, 2).Value = 2018
$ws.Cells.Item(1332, 3).Value = "Mainline (G Case)"
$ws.Cells.Item(1332, 4).Value = "Pontiac"
$ws.Cells.Item(1332, 5).Value = "Firebird Custom"
$ws.Cells.Item(1332, 6).Value = 1968
$ws.Cells.Item(1332, 7).Value = "Red"
$ws.Cells.Item(1332, 8).Value = "No"

# Row 1333
$ws.Cells.Item(1333, 1).Value = "HotWheels"
$ws.Cells.Item(1333, 2).Value = 2018
$ws.Cells.Item(1333, 3).Value = "Mainline (G Case)"
$ws.Cells.Item(1333, 4).Value = "Pontiac"
$ws.Cells.Item(1333, 5).Value = "Firebird Custom"
$ws.Cells.Item(1333, 6).Value = 1968
$ws.Cells.Item(1333, 7).Value = "Red"
$ws.Cells.Item(1333, 8).Value = "Yes"

# Row 1334
$ws.Cells.Item(1334, 1).Value = "HotWheels"
$ws.Cells.Item(1334, 2).Value = 2017
$ws.Cells.Item(1334, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1334, 4).Value = "Ford"
$ws.Cells.Item(1334, 5).Value = "Coupe"
$ws.Cells.Item(1334, 6).Value = 1932
$ws.Cells.Item(1334, 7).Value = "Silver"
$ws.Cells.Item(1334, 8).Value = "No"

# Row 1335
$ws.Cells.Item(1335, 1).Value = "HotWheels"
$ws.Cells.Item(1335, 2).Value = 2018
$ws.Cells.Item(1335, 3).Value = "Mainline (A Case)"
$ws.Cells.Item(1335, 4).Value = "Datsun"
$ws.Cells.Item(1335, 5).Value = "Fairlady 2000"
$ws.Cells.Item(1335, 6).Value = 1969
$ws.Cells.Item(1335, 7).Value = "Black"
$ws.Cells.Item(1335, 8).Value = "No"

# Row 1336
$ws.Cells.Item(1336, 1).Value = "HotWheels"
$ws.Cells.Item(1336, 2).Value = 2019
$ws.Cells.Item(1336, 3).Value = "Mainline (B Case)"
$ws.Cells.Item(1336, 4).Value = "Chevy"
$ws.Cells.Item(1336, 5).Value = "Custom Luv"
$ws.Cells.Item(1336, 6).Value = 1972
$ws.Cells.Item(1336, 7).Value = "Lime Green"
$ws.Cells.Item(1336, 8).Value = "No"

# Row 1337
$ws.Cells.Item(1337, 1).Value = "HotWheels"
$ws.Cells.Item(1337, 2).Value = 2017
$ws.Cells.Item(1337, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1337, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1337, 5).Value = "Torque Twister"
$ws.Cells.Item(1337, 6).Value = 2017
$ws.Cells.Item(1337, 7).Value = "White"
$ws.Cells.Item(1337, 8).Value = "No"

# Row 1338
$ws.Cells.Item(1338, 1).Value = "HotWheels"
$ws.Cells.Item(1338, 2).Value = 2017
$ws.Cells.Item(1338, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1338, 4).Value = "Ford"
$ws.Cells.Item(1338, 5).Value = "Coupe"
$ws.Cells.Item(1338, 6).Value = 1940
$ws.Cells.Item(1338, 7).Value = "Blue"
$ws.Cells.Item(1338, 8).Value = "No"

# Row 1339
$ws.Cells.Item(1339, 1).Value = "HotWheels"
$ws.Cells.Item(1339, 2).Value = 2017
$ws.Cells.Item(1339, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1339, 4).Value = "Plymouth"
$ws.Cells.Item(1339, 5).Value = "Fury"
$ws.Cells.Item(1339, 6).Value = 1957
$ws.Cells.Item(1339, 7).Value = "Lime Green"
$ws.Cells.Item(1339, 8).Value = "No"

# Row 1340
$ws.Cells.Item(1340, 1).Value = "HotWheels"
$ws.Cells.Item(1340, 2).Value = 2017
$ws.Cells.Item(1340, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1340, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1340, 5).Value = "Covelight"
$ws.Cells.Item(1340, 6).Value = 2017
$ws.Cells.Item(1340, 7).Value = "Red"
$ws.Cells.Item(1340, 8).Value = "No"

# Row 1341
$ws.Cells.Item(1341, 1).Value = "HotWheels"
$ws.Cells.Item(1341, 2).Value = 2017
$ws.Cells.Item(1341, 3).Value = "Mickey Mouse"
$ws.Cells.Item(1341, 4).Value = "Fake HotWheels Brand"
$ws.Cells.Item(1341, 5).Value = "Rocket Box"
$ws.Cells.Item(1341, 6).Value = 2017
$ws.Cells.Item(1341, 7).Value = "Gold"
$ws.Cells.Item(1341, 8).Value = "No"

# --- Update viewport selection to match the post-edit cursor position ---
$ws.Range("A1342").Select()
